$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H4").Value = "Test Value"
$ws.Range("H4").AddComment("Test Comment") | Out-Null
